$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Finish the literal CSV import: the bracketed gloss "[Hev, Ro]" should be
# split into two separate bracketed tags "[Hev] [Ro]" in the two cells
# that still contained the old-style combined tag (B4 and C4). Both cells
# are rich-text (mixed formatting runs), so the substitution is performed
# with Characters() so any existing run formatting on the untouched text
# is left alone.

$oldTag = "[Hev, Ro]"
$newTag = "[Hev] [Ro]"

$b4 = $ws.Range("B4")
$b4Text = $b4.Characters().Text
$b4Start = $b4Text.IndexOf($oldTag) + 1
if ($b4Start -gt 0) {
    $b4.Characters($b4Start, $oldTag.Length).Text = $newTag
}

$c4 = $ws.Range("C4")
$c4Text = $c4.Characters().Text
$c4Start = $c4Text.IndexOf($oldTag) + 1
if ($c4Start -gt 0) {
    $c4.Characters($c4Start, $oldTag.Length).Text = $newTag
}

# Reset the sheet's selection back to the default top-left cell (the
# import previously left the cursor parked at D12).
$ws.Range("A1").Select()
